$wb = $excel.ActiveWorkbook

# Sheets "展览" (Exhibitions) and "全部类型" (All types) both hold the same
# event table and both receive identical updates in this diff.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Update "想去人数" (interest count) figures in column F ---
    $ws.Range("F2").Value2  = 8374
    $ws.Range("F3").Value2  = 7915
    $ws.Range("F5").Value2  = 192
    $ws.Range("F10").Value2 = 177
    $ws.Range("F11").Value2 = 232
    $ws.Range("F13").Value2 = 134
    $ws.Range("F14").Value2 = 1899
    $ws.Range("F16").Value2 = 57

    # --- Append a new event row (row 20) mirroring the formatting of row 19 ---
    $ws.Range("A19:I19").Copy()
    $ws.Range("A20:I20").PasteSpecial(-4122)

    $ws.Range("A20").Value2 = 19

    $ws.Range("B20").Value2 = "'2024-06-08"
    $ws.Range("C20").Value2 = "'合肥·环形宇宙动漫游戏嘉年华-一周年超强巨制~"
    $ws.Range("D20").Value2 = "'锦绣大道3899号 合肥滨湖会展中心"
    $ws.Range("E20").Value2 = "'2024.06.08 09:30-06.09 17:00"
    $ws.Range("F20").Value2 = 2
    $ws.Range("G20").Value2 = 65
    $ws.Range("H20").Value2 = "'https://show.bilibili.com/platform/detail.html?id=83518"
    $ws.Range("I20").Value2 = "'//i1.hdslb.com/bfs/openplatform/202403/1Sqp42gM1711691520194.jpeg"

    # Reset style on the text cells so the apostrophe-prefix "quote prefix"
    # formatting doesn't stick (matches plain default style of the source rows).
    $ws.Range("B20:E20").Style = "Normal"
    $ws.Range("H20:I20").Style = "Normal"
}
